$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 64
    $ws.Range("F7").Value = 31
    $ws.Range("F9").Value = 49
    $ws.Range("F10").Value = 3
    $ws.Range("F11").Value = 4608
    $ws.Range("F12").Value = 4421
}
